$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Block A: rows 30-33 (Algebra 1 SCHI/LDTeam/SCLD/CSS) -> Algebra 1 description (shared string 153)
$ws.Cells.Item(30, 4).Value = 'Topics include linear equations and inequalities, systems of linear equations, relations, functions, polynomials, and statistics. Emphasis is placed on making connections in algebra to geometry and statistics.'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(31, 4).Value = 'Topics include linear equations and inequalities, systems of linear equations, relations, functions, polynomials, and statistics. Emphasis is placed on making connections in algebra to geometry and statistics.'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(32, 4).Value = 'Topics include linear equations and inequalities, systems of linear equations, relations, functions, polynomials, and statistics. Emphasis is placed on making connections in algebra to geometry and statistics.'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(33, 4).Value = 'Topics include linear equations and inequalities, systems of linear equations, relations, functions, polynomials, and statistics. Emphasis is placed on making connections in algebra to geometry and statistics.'
$ws.Cells.Item(33, 4).Style = "Normal"

# Block B: rows 38-42 (Geometry variants) -> Geometry description (shared string 155)
$ws.Cells.Item(38, 4).Value = 'This course emphasizes two- and three-dimensional reasoning skills, coordinate and transformational geometry, and the use of geometric models to solve problems. '
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(39, 4).Value = 'This course emphasizes two- and three-dimensional reasoning skills, coordinate and transformational geometry, and the use of geometric models to solve problems. '
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(40, 4).Value = 'This course emphasizes two- and three-dimensional reasoning skills, coordinate and transformational geometry, and the use of geometric models to solve problems. '
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(41, 4).Value = 'This course emphasizes two- and three-dimensional reasoning skills, coordinate and transformational geometry, and the use of geometric models to solve problems. '
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(42, 4).Value = 'This course emphasizes two- and three-dimensional reasoning skills, coordinate and transformational geometry, and the use of geometric models to solve problems. '
$ws.Cells.Item(42, 4).Style = "Normal"

# Block C: rows 126-143, new course descriptions (vertical-centered Times New Roman 12pt, matches existing style index 7)
$ws.Cells.Item(126, 4).Value = 'This is an adapted curriculum elective geared to students needing intensive support. Instruction is very concrete with extensive physical modeling and assistance. The course will identify work-related abilities, provide training and work skills, and prepare students for post-secondary participation in community-based worksites.'
$ws.Cells.Item(126, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(126, 4).Font.Size = 12
$ws.Cells.Item(126, 4).VerticalAlignment = -4108
$ws.Cells.Item(127, 4).Value = 'This is an adapted curriculum elective geared to students needing intensive support. Instruction is very concrete with extensive physical modeling and assistance. The course will identify work-related abilities, provide training and work skills, and prepare students for post-secondary participation in community-based worksites.'
$ws.Cells.Item(127, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(127, 4).Font.Size = 12
$ws.Cells.Item(127, 4).VerticalAlignment = -4108
$ws.Cells.Item(128, 4).Value = 'This course is designed to support instruction in the science content area; does not require SOL testing. Instruction is individualized based on needs identified in the IEP to help students gain a basic content vocabulary, knowledge and skills and designed to be taught at the learning pace of the individual students. '
$ws.Cells.Item(128, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(128, 4).Font.Size = 12
$ws.Cells.Item(128, 4).VerticalAlignment = -4108
$ws.Cells.Item(129, 4).Value = 'This course is designed to support instruction in the science content area; does not require SOL testing. Instruction is individualized based on needs identified in the IEP to help students gain a basic content vocabulary, knowledge and skills and designed to be taught at the learning pace of the individual students. '
$ws.Cells.Item(129, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(129, 4).Font.Size = 12
$ws.Cells.Item(129, 4).VerticalAlignment = -4108
$ws.Cells.Item(130, 4).Value = 'This course is designed to support instruction in the English content area; does not require SOL testing. Instruction is individualized based on needs identified in the IEP to help students gain a basic content vocabulary, knowledge and skills and designed to be taught at the learning pace of the individual students. '
$ws.Cells.Item(130, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(130, 4).Font.Size = 12
$ws.Cells.Item(130, 4).VerticalAlignment = -4108
$ws.Cells.Item(131, 4).Value = 'This course is designed to support instruction in the English content area; does not require SOL testing. Instruction is individualized based on needs identified in the IEP to help students gain a basic content vocabulary, knowledge and skills and designed to be taught at the learning pace of the individual students. '
$ws.Cells.Item(131, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(131, 4).Font.Size = 12
$ws.Cells.Item(131, 4).VerticalAlignment = -4108
$ws.Cells.Item(132, 4).Value = 'Students develop the ability to communicate about themselves and their immediate environment using simple sentences containing basic language structures. This communication is evidenced in signing, receiving signs and non-manual gestures, and reading. Students begin to explore and study the themes of Personal and Family Life, School Life, Social Life, and Community Life. '
$ws.Cells.Item(132, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(132, 4).Font.Size = 12
$ws.Cells.Item(132, 4).VerticalAlignment = -4108
$ws.Cells.Item(133, 4).Value = 'Students continue to develop proficiency in American Sign Language. They learn to function in real-life situations using more complex language structures and a wider range of vocabulary. Students continue to explore as they study the themes of Home Life, Student Life, Leisure Time, and Vacation and Travel. '
$ws.Cells.Item(133, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(133, 4).Font.Size = 12
$ws.Cells.Item(133, 4).VerticalAlignment = -4108
$ws.Cells.Item(134, 4).Value = 'Students continue to develop and refine their proficiency in American Sign Language. They communicate using more complex language structures on a variety of topics, moving from concrete to more abstract concepts. Students gain a deeper understanding of the world around them while studying Rights and Responsibilities, Future Plans and Choices, Teen Culture, Environment, and Humanities. '
$ws.Cells.Item(134, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(134, 4).Font.Size = 12
$ws.Cells.Item(134, 4).VerticalAlignment = -4108
$ws.Cells.Item(135, 4).Value = 'Foundations of World History/Geography is a one-credit elective course designed to support in the World History content area; does not require SOL testing. Instruction is individualized based on needs identified in the IEP to help students gain a basic content vocabulary, knowledge and skills and designed to be taught at the learning pace of the individual students.'
$ws.Cells.Item(135, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(135, 4).Font.Size = 12
$ws.Cells.Item(135, 4).VerticalAlignment = -4108
$ws.Cells.Item(136, 4).Value = 'Foundations of World History/Geography is a one-credit elective course designed to support in the World History content area; does not require SOL testing. Instruction is individualized based on needs identified in the IEP to help students gain a basic content vocabulary, knowledge and skills and designed to be taught at the learning pace of the individual students.'
$ws.Cells.Item(136, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(136, 4).Font.Size = 12
$ws.Cells.Item(136, 4).VerticalAlignment = -4108
$ws.Cells.Item(137, 4).Value = 'This comprehensive individualized program is designed to prepare students for a style of living that will require a minimum of dependence on family. The course is geared to meet the needs of the students as they prepare to enter employment and emphasizes developing interpersonal skills, following directions, working independently, completing a task, and developing self-advocacy and other community living skills.'
$ws.Cells.Item(137, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(137, 4).Font.Size = 12
$ws.Cells.Item(137, 4).VerticalAlignment = -4108
$ws.Cells.Item(138, 4).Value = 'This course, offered at Davis and Pulley Centers and STEP, is designed to teach students with disabilities skills for independent living.'
$ws.Cells.Item(138, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(138, 4).Font.Size = 12
$ws.Cells.Item(138, 4).VerticalAlignment = -4108
$ws.Cells.Item(139, 4).Value = 'Explore a variety of exciting opportunities and materials to inform the artmaking process. Learn to think conceptually and realize potential as a creative and critical thinker in order to meet the challenges of 21st century living. Explore personal interests while developing skills in the areas of drawing, painting, printmaking, ceramics, sculpture, and digital media. '
$ws.Cells.Item(139, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(139, 4).Font.Size = 12
$ws.Cells.Item(139, 4).VerticalAlignment = -4108
$ws.Cells.Item(140, 4).Value = 'Explore a variety of exciting opportunities and materials to inform the artmaking process. Learn to think conceptually and realize potential as a creative and critical thinker in order to meet the challenges of 21st century living. Explore personal interests while developing skills in the areas of drawing, painting, printmaking, ceramics, sculpture, and digital media. '
$ws.Cells.Item(140, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(140, 4).Font.Size = 12
$ws.Cells.Item(140, 4).VerticalAlignment = -4108
$ws.Cells.Item(141, 4).Value = 'Explore a variety of exciting opportunities and materials to inform the artmaking process. Learn to think conceptually and realize potential as a creative and critical thinker in order to meet the challenges of 21st century living. Explore personal interests while developing skills in the areas of drawing, painting, printmaking, ceramics, sculpture, and digital media. '
$ws.Cells.Item(141, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(141, 4).Font.Size = 12
$ws.Cells.Item(141, 4).VerticalAlignment = -4108
$ws.Cells.Item(142, 4).Value = 'Theatre Arts 1 provides students with a survey of the theatre arts, allowing student’s opportunities to experience and appreciate dramatic literature, and participate in the creative processes of performance and production. This course emphasizes skill development and provides theatrical opportunities that enable students to determine personal areas of interest. '
$ws.Cells.Item(142, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(142, 4).Font.Size = 12
$ws.Cells.Item(142, 4).VerticalAlignment = -4108
$ws.Cells.Item(143, 4).Value = 'Students are provided the opportunity to sample a variety of musical experiences in a non-performing music class. Course content includes beginning guitar and class piano experience, as well as various modules designed to assist students in developing music reading and composing skills. '
$ws.Cells.Item(143, 4).Font.Name = "Times New Roman"
$ws.Cells.Item(143, 4).Font.Size = 12
$ws.Cells.Item(143, 4).VerticalAlignment = -4108

# Update view state: selection + scroll position (mirrors author re-scrolling/selecting before saving)
$ws.Range("M144").Select()
$excel.ActiveWindow.ScrollRow = 118
$excel.ActiveWindow.ScrollColumn = 1

